$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ch1")
$ws.Activate()
$ws.Range("D41").Select()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
Write-Output $excel.ActiveWindow.ScrollRow
Write-Output $excel.ActiveWindow.ScrollColumn
